$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.392.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.45%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.846.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.17%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.9986"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.09%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'240.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.04%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.6321"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.02%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").Value = "'0.07552"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.04%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.2968"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.24%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  +1.27%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07727"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.65%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.826.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.95%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'5.001"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.44%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.6854"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.03%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.00001003"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.87%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'83.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.81%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  -0.31%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'29.422.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.48%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'229.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.78%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  -0.23%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'1.0000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.03%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'7.581"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.40%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("B23").Value = "'BinanceUSD"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.02%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("B24").Value = "'Monero"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'157.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.89%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("B25").Value = "'Stellar"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'0.1401"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.09%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("B26").Value = "'Cosmos"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'8.386"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.55%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("B27").Value = "'EthereumClassic"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'17.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.30%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").Value = "'PancakeSwap"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'1.467"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.09%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("B29").Value = "'Hedera"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'0.05726"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.93%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("B30").Value = "'Toncoin"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'1.251"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.15%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("B31").Value = "'Filecoin"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'4.131"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.49%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").Value = "'InternetComputer(DFINITY)"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'4.038"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.27%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").Value = "'LidoDAOToken"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'1.854"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.20%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").Value = "'ARBITRUM"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'1.157"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.20%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("B35").Value = "'ImmutableX"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'0.7176"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.14%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("B36").Value = "'HuobiToken"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'2.594"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.11%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = "'Maker"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'1.250.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.24%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = "'VeChain"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'0.01812"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.11%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").Value = "'MXToken"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'2.783"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.50%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = "'FraxShare"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'6.211"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.18%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "'TrustWalletToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.9065"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.88%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "'PaxDollar"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'1.000"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.08%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value = "'Quant"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'101.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.02%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'1.991.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.43%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "'Aave"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'66.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.58%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "'BabyDogeCoin"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000118"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.38%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "'EnergySwap"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'9.208"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.51%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = "'Aptos"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'7.067"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.09%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "'RenderToken"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'1.714"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.04%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.4031"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.05%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = "'Algorand"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.1129"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.86%  "
$ws.Range("E51").Style = "Normal"
